$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(213, 2).Value = '['
$ws.Cells.Item(214, 2).Value = '    {'
$ws.Cells.Item(215, 2).Value = '        "name": "surname",'
$ws.Cells.Item(216, 2).Value = '        "value": "MONTEALEGRE"'
$ws.Cells.Item(217, 2).Value = '    },'
$ws.Cells.Item(218, 2).Value = '    {'
$ws.Cells.Item(219, 2).Value = '        "name": "firstName",'
$ws.Cells.Item(220, 2).Value = '        "value": "ASUNCION"'
$ws.Cells.Item(221, 2).Value = '    },'
$ws.Cells.Item(222, 2).Value = '    {'
$ws.Cells.Item(223, 2).Value = '        "name": "middlename",'
$ws.Cells.Item(224, 2).Value = '        "value": "SABANDAL"'
$ws.Cells.Item(225, 2).Value = '    },'
$ws.Cells.Item(226, 2).Value = '    {'
$ws.Cells.Item(227, 2).Value = '        "name": "barangay",'
$ws.Cells.Item(228, 2).Value = '        "value": "Tagapo"'
$ws.Cells.Item(229, 2).Value = '    },'
$ws.Cells.Item(230, 2).Value = '    {'
$ws.Cells.Item(231, 2).Value = '        "name": "address",'
$ws.Cells.Item(232, 2).Value = '        "value": "BLOCK 7 LOT 2 OAK STREET ROSE POINTE SUBDIVISION"'
$ws.Cells.Item(233, 2).Value = '    },'
$ws.Cells.Item(234, 2).Value = '    {'
$ws.Cells.Item(235, 2).Value = '        "name": "gender",'
$ws.Cells.Item(236, 2).Value = '        "value": "Female"'
$ws.Cells.Item(237, 2).Value = '    },'
$ws.Cells.Item(238, 2).Value = '    {'
$ws.Cells.Item(239, 2).Value = '        "name": "maritalStatus",'
$ws.Cells.Item(240, 2).Value = '        "value": "Widowed"'
$ws.Cells.Item(241, 2).Value = '    },'
$ws.Cells.Item(242, 2).Value = '    {'
$ws.Cells.Item(243, 2).Value = '        "name": "srCitizenDOB",'
$ws.Cells.Item(244, 2).Value = '        "value": "1942-01-12"'
$ws.Cells.Item(245, 2).Value = '    },'
$ws.Cells.Item(246, 2).Value = '    {'
$ws.Cells.Item(247, 2).Value = '        "name": "placeOfBirth",'
$ws.Cells.Item(248, 2).Value = '        "value": "AKLAN"'
$ws.Cells.Item(249, 2).Value = '    },'
$ws.Cells.Item(250, 2).Value = '    {'
$ws.Cells.Item(251, 2).Value = '        "name": "email",'
$ws.Cells.Item(252, 2).Value = '        "value": "asuncionmontealegre@gmail.com"'
$ws.Cells.Item(253, 2).Value = '    },'
$ws.Cells.Item(254, 2).Value = '    {'
$ws.Cells.Item(255, 2).Value = '        "name": "telephone",'
$ws.Cells.Item(256, 2).Value = '        "value": "09091072865"'
$ws.Cells.Item(257, 2).Value = '    },'
$ws.Cells.Item(258, 2).Value = '    {'
$ws.Cells.Item(259, 2).Value = '        "name": "religion",'
$ws.Cells.Item(260, 2).Value = '        "value": "CATHOLIC"'
$ws.Cells.Item(261, 2).Value = '    },'
$ws.Cells.Item(262, 2).Value = '    {'
$ws.Cells.Item(263, 2).Value = '        "name": "job",'
$ws.Cells.Item(264, 2).Value = '        "value": "NONE"'
$ws.Cells.Item(265, 2).Value = '    },'
$ws.Cells.Item(266, 2).Value = '    {'
$ws.Cells.Item(267, 2).Value = '        "name": "hasPension",'
$ws.Cells.Item(268, 2).Value = '        "value": "Meron"'
$ws.Cells.Item(269, 2).Value = '    },'
$ws.Cells.Item(270, 2).Value = '    {'
$ws.Cells.Item(271, 2).Value = '        "name": "whatPension",'
$ws.Cells.Item(272, 2).Value = '        "value": "SSS"'
$ws.Cells.Item(273, 2).Value = '    },'
$ws.Cells.Item(274, 2).Value = '    {'
$ws.Cells.Item(275, 2).Value = '        "name": "howMuchPension",'
$ws.Cells.Item(276, 2).Value = '        "value": "6000"'
$ws.Cells.Item(277, 2).Value = '    },'
$ws.Cells.Item(278, 2).Value = '    {'
$ws.Cells.Item(279, 2).Value = '        "name": "spouseLastName",'
$ws.Cells.Item(280, 2).Value = '        "value": "MONTEALEGRE"'
$ws.Cells.Item(281, 2).Value = '    },'
$ws.Cells.Item(282, 2).Value = '    {'
$ws.Cells.Item(283, 2).Value = '        "name": "spouseFirstName",'
$ws.Cells.Item(284, 2).Value = '        "value": "VERGILIO"'
$ws.Cells.Item(285, 2).Value = '    },'
$ws.Cells.Item(286, 2).Value = '    {'
$ws.Cells.Item(287, 2).Value = '        "name": "spouseMiddleName",'
$ws.Cells.Item(288, 2).Value = '        "value": "SOBERANO"'
$ws.Cells.Item(289, 2).Value = '    },'
$ws.Cells.Item(290, 2).Value = '    {'
$ws.Cells.Item(291, 2).Value = '        "name": "spouseSuffix",'
$ws.Cells.Item(292, 2).Value = '        "value": ""'
$ws.Cells.Item(293, 2).Value = '    },'
$ws.Cells.Item(294, 2).Value = '    {'
$ws.Cells.Item(295, 2).Value = '        "name": "spouseDOB",'
$ws.Cells.Item(296, 2).Value = '        "value": "1942-08-07"'
$ws.Cells.Item(297, 2).Value = '    },'
$ws.Cells.Item(298, 2).Value = '    {'
$ws.Cells.Item(299, 2).Value = '        "name": "numberOfChildren",'
$ws.Cells.Item(300, 2).Value = '        "value": "1"'
$ws.Cells.Item(301, 2).Value = '    },'
$ws.Cells.Item(302, 2).Value = '    {'
$ws.Cells.Item(303, 2).Value = '        "name": "totalHousemate",'
$ws.Cells.Item(304, 2).Value = '        "value": "1"'
$ws.Cells.Item(305, 2).Value = '    },'
$ws.Cells.Item(306, 2).Value = '    {'
$ws.Cells.Item(307, 2).Value = '        "name": "childFirstName",'
$ws.Cells.Item(308, 2).Value = '        "value": "MONTEALEGRE"'
$ws.Cells.Item(309, 2).Value = '    },'
$ws.Cells.Item(310, 2).Value = '    {'
$ws.Cells.Item(311, 2).Value = '        "name": "childLastName",'
$ws.Cells.Item(312, 2).Value = '        "value": "MILA ROSA"'
$ws.Cells.Item(313, 2).Value = '    },'
$ws.Cells.Item(314, 2).Value = '    {'
$ws.Cells.Item(315, 2).Value = '        "name": "srCitizenChildDOB",'
$ws.Cells.Item(316, 2).Value = '        "value": "1965-07-20"'
$ws.Cells.Item(317, 2).Value = '    },'
$ws.Cells.Item(318, 2).Value = '    {'
$ws.Cells.Item(319, 2).Value = '        "name": "childTelephone",'
$ws.Cells.Item(320, 2).Value = '        "value": "09760657071"'
$ws.Cells.Item(321, 2).Value = '    },'
$ws.Cells.Item(322, 2).Value = '    {'
$ws.Cells.Item(323, 2).Value = '        "name": "childBarangay",'
$ws.Cells.Item(324, 2).Value = '        "value": "Tagapo"'
$ws.Cells.Item(325, 2).Value = '    },'
$ws.Cells.Item(326, 2).Value = '    {'
$ws.Cells.Item(327, 2).Value = '        "name": "childAddress",'
$ws.Cells.Item(328, 2).Value = '        "value": "BLOCK 7 LOT 2 ROSE POINTE SUBD BRGY. TAGAPO SANTA ROSA LAGUNA"'
$ws.Cells.Item(329, 2).Value = '    }'
$ws.Cells.Item(330, 2).Value = ']'

[void]$ws.Range("B213:B330").Select()
$excel.ActiveWindow.ScrollRow = 199
